# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (Total) sheet,
#    with the quarterly fund-holdings table for 2022-Q1.
# 2. Prepend a "2022-Q1" row to the "总计" summary sheet and renumber the
#    existing index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New sheet "2022-Q1", inserted immediately before "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$headerRange = $q1.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data row (row 2) - numeric-looking text columns stay text, like the
# other quarter sheets (fund code / size / position figures are strings).
$q1.Range("A2").Value = 0
$q1.Range("A2").Font.Bold = $true
$q1.Range("A2").HorizontalAlignment = -4108
$q1.Range("A2").VerticalAlignment = -4160
$q1.Range("A2").Borders.LineStyle = 1

$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "519097"
$q1.Range("C2").Value = "新华中小市值优选混合"
$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "0.75"
$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "62.70"
$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "4.67"
$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "0.0350"
$q1.Range("H2").Value = 4

# ---------------------------------------------------------------------
# 2) "总计" sheet: add the 2022-Q1 row on top, shift the rest down
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()
# Insert() copies the header row's formatting down onto the blank row;
# strip it from the plain data cells (only column A carries the bold/
# bordered "index" look, same as every other row in this sheet).
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.04
$total.Range("A2").Font.Bold = $true
$total.Range("A2").HorizontalAlignment = -4108
$total.Range("A2").VerticalAlignment = -4160
$total.Range("A2").Borders.LineStyle = 1

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.31

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 17
$total.Range("D4").Value = 5.21

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 7
$total.Range("D5").Value = 1.08

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 13
$total.Range("D6").Value = 1.95

$total.Range("A7").Value = 5
$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 16
$total.Range("D7").Value = 2.52

# ---------------------------------------------------------------------
# Restore the original active sheet/tab selection
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
